# Updated cryptos list (price / 1h-volume refresh) matching the commit's
# OOXML diff. Column D/E cells whose new value reads as a plain number
# (e.g. "7.63") would otherwise be auto-converted to a numeric cell by
# Excel's normal type inference; since the source data stores these as
# text, such cells are first forced to text format ("@"), assigned, and
# then have ClearFormats() applied so the cell reverts to the sheet's
# default (unstyled) cell format while keeping its value stored as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '67.549.82'

$ws.Cells.Item(2, 5).Value = '  -0.08%  '

$ws.Cells.Item(3, 4).Value = '3.510.29'

$ws.Cells.Item(3, 5).Value = '  -0.43%  '

$ws.Cells.Item(4, 5).Value = '  +0.11%  '

$c = $ws.Cells.Item(5, 4)
$c.NumberFormat = "@"
$c.Value = '609.03'
$c.ClearFormats()

$ws.Cells.Item(5, 5).Value = '  -0.29%  '

$c = $ws.Cells.Item(6, 4)
$c.NumberFormat = "@"
$c.Value = '152.41'
$c.ClearFormats()

$ws.Cells.Item(6, 5).Value = '  +0.41%  '

$ws.Cells.Item(7, 4).Value = '3.509.31'

$ws.Cells.Item(7, 5).Value = '  -0.40%  '

$ws.Cells.Item(8, 5).Value = '  +0.00%  '

$ws.Cells.Item(9, 5).Value = '  +1.01%  '

$c = $ws.Cells.Item(10, 4)
$c.NumberFormat = "@"
$c.Value = '0.144'
$c.ClearFormats()

$ws.Cells.Item(10, 5).Value = '  +2.48%  '

$c = $ws.Cells.Item(11, 4)
$c.NumberFormat = "@"
$c.Value = '7.63'
$c.ClearFormats()

$ws.Cells.Item(11, 5).Value = '  +8.05%  '

$c = $ws.Cells.Item(12, 4)
$c.NumberFormat = "@"
$c.Value = '0.432'
$c.ClearFormats()

$ws.Cells.Item(12, 5).Value = '  +1.52%  '

$ws.Cells.Item(13, 5).Value = '  +2.24%  '

$ws.Cells.Item(14, 5).Value = '  -2.14%  '

$ws.Cells.Item(15, 4).Value = '4.104.52'

$ws.Cells.Item(15, 5).Value = '  -0.32%  '

$ws.Cells.Item(16, 4).Value = '3.509.89'

$ws.Cells.Item(16, 5).Value = '  -0.42%  '

$ws.Cells.Item(17, 4).Value = '67.542.52'

$ws.Cells.Item(17, 5).Value = '  +0.07%  '

$ws.Cells.Item(18, 5).Value = '  +0.42%  '

$ws.Cells.Item(19, 5).Value = '  +2.28%  '

$c = $ws.Cells.Item(20, 4)
$c.NumberFormat = "@"
$c.Value = '15.59'
$c.ClearFormats()

$ws.Cells.Item(20, 5).Value = '  +2.51%  '

$c = $ws.Cells.Item(21, 4)
$c.NumberFormat = "@"
$c.Value = '9.87'
$c.ClearFormats()

$ws.Cells.Item(21, 5).Value = '  +6.10%  '

$c = $ws.Cells.Item(22, 4)
$c.NumberFormat = "@"
$c.Value = '448.79'
$c.ClearFormats()

$ws.Cells.Item(22, 5).Value = '  +0.45%  '

$ws.Cells.Item(23, 5).Value = '  +1.50%  '

$c = $ws.Cells.Item(24, 4)
$c.NumberFormat = "@"
$c.Value = '78.22'
$c.ClearFormats()

$ws.Cells.Item(24, 5).Value = '  +1.05%  '

$ws.Cells.Item(25, 2).Value = 'WrappedeETH'

$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'

$ws.Cells.Item(25, 4).Value = '3.651.45'

$ws.Cells.Item(25, 5).Value = '  -0.41%  '

$ws.Cells.Item(26, 2).Value = 'PEPE'

$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'

$c = $ws.Cells.Item(26, 4)
$c.NumberFormat = "@"
$c.Value = '0.0000128'
$c.ClearFormats()

$ws.Cells.Item(26, 5).Value = '  -1.38%  '

$ws.Cells.Item(27, 5).Value = '  +0.00%  '

$c = $ws.Cells.Item(28, 4)
$c.NumberFormat = "@"
$c.Value = '8.86'
$c.ClearFormats()

$ws.Cells.Item(28, 5).Value = '  +5.77%  '

$c = $ws.Cells.Item(29, 4)
$c.NumberFormat = "@"
$c.Value = '10.12'
$c.ClearFormats()

$ws.Cells.Item(29, 5).Value = '  -0.90%  '

$ws.Cells.Item(30, 5).Value = '  +0.69%  '

$ws.Cells.Item(31, 5).Value = '  +5.77%  '

$c = $ws.Cells.Item(32, 4)
$c.NumberFormat = "@"
$c.Value = '0.168'
$c.ClearFormats()

$ws.Cells.Item(32, 5).Value = '  +2.17%  '

$ws.Cells.Item(33, 5).Value = '  +0.13%  '

$c = $ws.Cells.Item(34, 4)
$c.NumberFormat = "@"
$c.Value = '25.79'
$c.ClearFormats()

$ws.Cells.Item(34, 5).Value = '  -0.04%  '

$ws.Cells.Item(35, 5).Value = '  +0.94%  '

$ws.Cells.Item(36, 5).Value = '  +1.60%  '

$ws.Cells.Item(37, 4).Value = '3.503.47'

$ws.Cells.Item(37, 5).Value = '  -0.39%  '

$c = $ws.Cells.Item(38, 4)
$c.NumberFormat = "@"
$c.Value = '8.07'
$c.ClearFormats()

$ws.Cells.Item(38, 5).Value = '  +0.05%  '

$ws.Cells.Item(39, 5).Value = '  +0.02%  '

$c = $ws.Cells.Item(40, 4)
$c.NumberFormat = "@"
$c.Value = '2.30'
$c.ClearFormats()

$ws.Cells.Item(40, 5).Value = '  +4.72%  '

$ws.Cells.Item(41, 5).Value = '  +0.21%  '

$c = $ws.Cells.Item(42, 4)
$c.NumberFormat = "@"
$c.Value = '0.0900'
$c.ClearFormats()

$ws.Cells.Item(42, 5).Value = '  +2.55%  '

$c = $ws.Cells.Item(43, 4)
$c.NumberFormat = "@"
$c.Value = '174.15'
$c.ClearFormats()

$ws.Cells.Item(43, 5).Value = '  -1.77%  '

$ws.Cells.Item(44, 5).Value = '  +0.88%  '

$c = $ws.Cells.Item(45, 4)
$c.NumberFormat = "@"
$c.Value = '30.43'
$c.ClearFormats()

$ws.Cells.Item(45, 5).Value = '  +11.90%  '

$c = $ws.Cells.Item(46, 4)
$c.NumberFormat = "@"
$c.Value = '0.883'
$c.ClearFormats()

$ws.Cells.Item(46, 5).Value = '  +0.23%  '

$c = $ws.Cells.Item(47, 4)
$c.NumberFormat = "@"
$c.Value = '46.69'
$c.ClearFormats()

$ws.Cells.Item(47, 5).Value = '  +2.42%  '

$ws.Cells.Item(48, 5).Value = '  +3.05%  '

$c = $ws.Cells.Item(49, 4)
$c.NumberFormat = "@"
$c.Value = '2.53'
$c.ClearFormats()

$ws.Cells.Item(49, 5).Value = '  -3.80%  '

$c = $ws.Cells.Item(50, 4)
$c.NumberFormat = "@"
$c.Value = '7.69'
$c.ClearFormats()

$ws.Cells.Item(50, 5).Value = '  +1.23%  '

$c = $ws.Cells.Item(51, 4)
$c.NumberFormat = "@"
$c.Value = '0.255'
$c.ClearFormats()

$ws.Cells.Item(51, 5).Value = '  +2.82%  '
